$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.401.83'
$ws.Range("E2").Value = '  +3.06%  '
$ws.Range("D3").Value = '2.657.52'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.84'
$ws.Range("E5").Value = '  +2.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.32'
$ws.Range("E6").Value = '  +4.26%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +1.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.119'
$ws.Range("E9").Value = '  +8.41%  '
$ws.Range("E10").Value = '  +4.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.88'
$ws.Range("E11").Value = '  +3.52%  '
$ws.Range("E12").Value = '  +1.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.27'
$ws.Range("E13").Value = '  +6.25%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000188'
$ws.Range("E14").Value = '  +20.44%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.133.04'
$ws.Range("E15").Value = '  +2.10%  '
$ws.Range("D16").Value = '65.302.50'
$ws.Range("E16").Value = '  +3.08%  '
$ws.Range("D17").Value = '2.767.75'
$ws.Range("E17").Value = '  +6.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.70'
$ws.Range("E18").Value = '  +2.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.86'
$ws.Range("E19").Value = '  +2.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '355.33'
$ws.Range("E20").Value = '  +2.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.33'
$ws.Range("E21").Value = '  +5.98%  '
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.35'
$ws.Range("E23").Value = '  +0.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.73'
$ws.Range("E24").Value = '  +2.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.59'
$ws.Range("E25").Value = '  +3.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.68'
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.30'
$ws.Range("E27").Value = '  +3.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.164'
$ws.Range("E28").Value = '  +2.13%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0964'
$ws.Range("E29").Value = '  +13.65%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '540.74'
$ws.Range("E31").Value = '  -3.42%  '
$ws.Range("E32").Value = '  +3.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.81'
$ws.Range("E33").Value = '  +3.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.79'
$ws.Range("E34").Value = '  +11.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.49'
$ws.Range("E35").Value = '  +5.05%  '
$ws.Range("E36").Value = '  +3.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.06'
$ws.Range("E37").Value = '  +6.62%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '165.65'
$ws.Range("E38").Value = '  -0.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '20.27'
$ws.Range("E39").Value = '  +3.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.996'
$ws.Range("E40").Value = '  -0.34%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '167.70'
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.22'
$ws.Range("E43").Value = '  +6.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.13'
$ws.Range("E44").Value = '  +4.37%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0612'
$ws.Range("E45").Value = '  +4.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.43'
$ws.Range("E46").Value = '  +5.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.25'
$ws.Range("E47").Value = '  +10.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.657'
$ws.Range("E48").Value = '  +4.09%  '
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0986'
$ws.Range("E50").Value = '  +2.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.69'
$ws.Range("E51").Value = '  +2.55%  '
